$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

# Update the A6 test script name
$ws.Range("A6").Value = "Dispensary\TC001CreateDispensarySale.py"

# Update the Env column (H) so every data row reads MPH
$ws.Range("H2:H18").Value = "MPH"

# Nudge column widths to match the resaved layout as closely as possible
$ws.Columns.Item(1).ColumnWidth = 53.666666666666664
$ws.Range("C1:D1").ColumnWidth = 10.5
$ws.Columns.Item(6).ColumnWidth = 16
$ws.Columns.Item(7).ColumnWidth = 12.833333333333334

# Update the current selection to mirror the saved workbook state
$ws.Range("A6").Select()
